# Commit: "make the wheel form canbe configed"
# The "幸运转盘" (lucky wheel) tip row is being removed from the LevelInfo
# table since the wheel's unlock text is now configurable elsewhere rather
# than hard-coded as a row in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (A8:E8) holds Id=5, Level=12, Type=1, Des="幸运转盘" tip, Icon=4.
# Deleting the entire row shifts rows 9-12 up, shrinks the worksheet
# dimension/table range by one row, and lets Excel prune the now-unused
# shared string for the wheel tip text.
$ws.Rows.Item(8).Delete()
